$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 170, shifting existing rows 170-175 down to 171-176
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new data record
$ws.Cells.Item(170, 1).Value()  = 4
$ws.Cells.Item(170, 2).Value()  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(170, 3).Value()  = "Los Lagos"
$ws.Cells.Item(170, 4).Value()  = 44706
$ws.Cells.Item(170, 5).Value()  = 10
$ws.Cells.Item(170, 6).Value()  = 100112009
$ws.Cells.Item(170, 7).Value()  = "Acelga"
$ws.Cells.Item(170, 8).Value()  = "Sin especificar"
$ws.Cells.Item(170, 9).Value()  = "Primera"
$ws.Cells.Item(170, 10).Value() = 20
$ws.Cells.Item(170, 11).Value() = 12000
$ws.Cells.Item(170, 12).Value() = 12000
$ws.Cells.Item(170, 13).Value() = 12000
$ws.Cells.Item(170, 14).Value() = "$/docena de atados (12 kilos)"
$ws.Cells.Item(170, 15).Value() = "Región de La Araucanía"
$ws.Cells.Item(170, 16).Value() = 1000
$ws.Cells.Item(170, 17).Value() = 12
$ws.Cells.Item(170, 18).Value() = "Hortaliza"
